# Update crypto price/volume values on the active sheet to reflect the latest
# scrape performed by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates that are NOT plain-number-looking strings (they
# already contain multiple '.' separators, e.g. "35.066.41") - these remain
# text automatically since Excel cannot parse them as a number.
$priceTextUpdates = @{
    2  = "35.066.41"
    3  = "1.858.16"
    12 = "2.127.01"
    13 = "1.865.30"
    17 = "35.043.47"
    41 = "1.354.14"
    49 = "2.040.19"
}

# Price (column D) updates that DO look like plain numbers (e.g. "238.36")
# - force the cell to Text format first so Excel doesn't convert the
# string into a numeric value, preserving the original text semantics.
$priceNumericLookingUpdates = @{
    5  = "238.36"
    8  = "42.49"
    14 = "11.41"
    15 = "0.678"
    16 = "4.71"
    18 = "70.29"
    20 = "241.04"
    24 = "2.27"
    25 = "171.45"
    27 = "7.92"
    28 = "17.72"
    36 = "0.790"
    39 = "91.71"
    42 = "14.91"
    43 = "2.35"
    44 = "12.85"
    47 = "0.0542"
    48 = "6.40"
}

# Volume(1h) (column E) updates - always text (percent strings with padding).
$volumeUpdates = @{
    2  = "  +1.60%  "
    3  = "  +3.27%  "
    4  = "  +0.17%  "
    5  = "  +4.06%  "
    7  = "  +0.16%  "
    8  = "  +8.64%  "
    9  = "  +3.30%  "
    10 = "  +2.99%  "
    11 = "  +0.32%  "
    12 = "  +3.25%  "
    13 = "  +3.95%  "
    15 = "  +3.21%  "
    16 = "  +3.58%  "
    17 = "  +1.97%  "
    18 = "  +1.98%  "
    19 = "  +2.64%  "
    20 = "  +0.65%  "
    21 = "  +3.31%  "
    22 = "  +1.45%  "
    24 = "  +2.14%  "
    25 = "  -0.56%  "
    26 = "  +31.19%  "
    27 = "  +3.07%  "
    28 = "  +3.29%  "
    29 = "  +2.89%  "
    30 = "  +3.41%  "
    31 = "  +0.15%  "
    32 = "  -0.10%  "
    33 = "  +3.52%  "
    34 = "  +13.75%  "
    35 = "  +22.58%  "
    36 = "  +14.39%  "
    37 = "  +6.01%  "
    38 = "  +12.90%  "
    39 = "  +1.21%  "
    40 = "  +6.63%  "
    41 = "  +3.20%  "
    42 = "  +4.55%  "
    43 = "  +6.73%  "
    44 = "  +54.80%  "
    45 = "  -0.12%  "
    46 = "  +1.76%  "
    47 = "  +5.71%  "
    48 = "  +4.09%  "
    49 = "  +2.91%  "
    50 = "  +3.30%  "
    51 = "  +18.08%  "
}

foreach ($row in $priceTextUpdates.Keys) {
    $ws.Range("D$row").Value = $priceTextUpdates[$row]
}

foreach ($row in $priceNumericLookingUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceNumericLookingUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
